$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the three new "GigabitEthernet1 / 10.10.20.x /24" rows.
#    Inserting shifts everything below down and copies the format of the row
#    above, which we will fix up explicitly afterwards.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(13).Insert()

# ---------------------------------------------------------------------------
# 2. Write the full final data grid (values only - formatting follows below).
# ---------------------------------------------------------------------------
$data = @(
    @(1,  "DISPOSITIVO", "INTERFAZ",         "IPv4",           "MÁSCARA", "GATEWAY"),
    @(2,  "R1",           "GigabitEthernet1", "10.10.20.181",   "/24",     ""),
    @(3,  "R1",           "GigabitEthernet2", "172.16.0.1",     "/30",     ""),
    @(4,  "R1",           "GigabitEthernet3", "172.16.0.9",     "/30",     ""),
    @(5,  "DISPOSITIVO", "INTERFAZ",         "IPv4",           "MÁSCARA", "GATEWAY"),
    @(6,  "R2",           "GigabitEthernet1", "10.10.20.182",   "/24",     ""),
    @(7,  "R2",           "GigabitEthernet2", "172.16.0.2",     "/30",     ""),
    @(8,  "R2",           "GigabitEthernet3", "172.16.0.5",     "/30",     ""),
    @(9,  "",             "",                 "",               "",        ""),
    @(10, "",             "",                 "",               "",        ""),
    @(11, "DISPOSITIVO", "INTERFAZ",         "IPv4",           "MÁSCARA", "GATEWAY"),
    @(12, "R3",           "GigabitEthernet1", "10.10.20.183",   "/24",     ""),
    @(13, "R3",           "GigabitEthernet2", "172.16.0.10",    "/30",     ""),
    @(14, "R3",           "GigabitEthernet3", "172.16.0.6",     "/30",     ""),
    @(15, "R3",           "GigabitEthernet4", "192.168.0.62",   "/26",     ""),
    @(16, "R3",           "GigabitEthernet5", "192.168.0.126",  "/26",     "")
)

foreach ($row in $data) {
    $r = $row[0]
    for ($c = 1; $c -le 5; $c++) {
        $val = $row[$c]
        $ws.Cells.Item($r, $c).Value = $val
    }
}

# ---------------------------------------------------------------------------
# 3. Formatting.
#
#    - Header rows (1, 5, 11) keep their original bold style (style index 1);
#      re-touch them so insert-artifacts are overwritten.
#    - "Normal" data rows get a thin border added all around (border index 1)
#      while keeping the existing centred alignment.
#    - The brand new "GigabitEthernet1" rows (2, 6, 12) additionally get a
#      light (theme background 1) fill, on top of the border + centring.
# ---------------------------------------------------------------------------

$newRows     = @(2, 6, 12)
$normalRows  = @(3, 4, 7, 8, 9, 10, 13, 14, 15, 16)

# The rows we just inserted inherited the bold header style from the row
# above. Reset them back to the plain body style first (copy format from an
# untouched body row) so subsequent fill/border tweaks build on top of the
# correct base font.
$plain = $ws.Range("A4:E4")
$plain.Copy()
foreach ($r in ($newRows + $normalRows)) {
    $ws.Range("A$r`:E$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

foreach ($r in $normalRows) {
    $rng = $ws.Range("A$r`:E$r")
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
}

foreach ($r in $newRows) {
    $rng = $ws.Range("A$r`:E$r")
    $rng.Interior.ThemeColor = 2
    $rng.Borders.LineStyle = 1
    $rng.HorizontalAlignment = -4108
}

# ---------------------------------------------------------------------------
# 4. Selection / dimension bookkeeping, matching the author's final state.
# ---------------------------------------------------------------------------
$ws.Range("C16").Select()
